# Apply updates described by the diff to worksheet "FTNT"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("C4").Value = 140000000.0
$ws.Range("D4").Value = 135000000.0
$ws.Range("E4").Value = 126000000.0
$ws.Range("F4").Value = 105000000.0
$ws.Range("G4").Value = 118000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 142000000.0
$ws.Range("D14").Value = 96000000.0
$ws.Range("E14").Value = 107000000.0
$ws.Range("F14").Value = 88000000.0
$ws.Range("G14").Value = 96000000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("C23").Value = -245000000.0
$ws.Range("D23").Value = -230000000.0
$ws.Range("E23").Value = -219000000.0
$ws.Range("F23").Value = -227000000.0
$ws.Range("G23").Value = -233000000.0

# Row 36 - Net Debt (B36 was an empty inline string, now a number)
$ws.Range("B36").Value = -1949900000.0

# Row 37 - Total Debt (B37 was an empty inline string, now a number)
$ws.Range("B37").Value = 987000000.0
